# Auto-generated edit script: shifts the "Fecha" (D) and "Origen" (O) values
# for rows 98-159 down the historical sequence by one reporting period, and
# appends two new rows (160-161) for Cilantro "Primera"/"Segunda" quality with
# the data point that used to be the most recent one (matching the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Fecha) updates for rows 98-159 ---
$ws.Cells.Item(98, 4).Value = 44582
$ws.Cells.Item(99, 4).Value = 44582
$ws.Cells.Item(100, 4).Value = 44308
$ws.Cells.Item(101, 4).Value = 44308
$ws.Cells.Item(102, 4).Value = 44320
$ws.Cells.Item(103, 4).Value = 44320
$ws.Cells.Item(104, 4).Value = 44306
$ws.Cells.Item(105, 4).Value = 44306
$ws.Cells.Item(106, 4).Value = 44316
$ws.Cells.Item(107, 4).Value = 44316
$ws.Cells.Item(108, 4).Value = 44460
$ws.Cells.Item(109, 4).Value = 44460
$ws.Cells.Item(110, 4).Value = 44467
$ws.Cells.Item(111, 4).Value = 44467
$ws.Cells.Item(112, 4).Value = 44313
$ws.Cells.Item(113, 4).Value = 44313
$ws.Cells.Item(114, 4).Value = 44334
$ws.Cells.Item(115, 4).Value = 44334
$ws.Cells.Item(116, 4).Value = 44209
$ws.Cells.Item(117, 4).Value = 44209
$ws.Cells.Item(118, 4).Value = 44405
$ws.Cells.Item(119, 4).Value = 44405
$ws.Cells.Item(120, 4).Value = 44280
$ws.Cells.Item(121, 4).Value = 44280
$ws.Cells.Item(122, 4).Value = 44330
$ws.Cells.Item(123, 4).Value = 44330
$ws.Cells.Item(124, 4).Value = 44239
$ws.Cells.Item(125, 4).Value = 44239
$ws.Cells.Item(126, 4).Value = 44476
$ws.Cells.Item(127, 4).Value = 44476
$ws.Cells.Item(128, 4).Value = 44169
$ws.Cells.Item(129, 4).Value = 44169
$ws.Cells.Item(130, 4).Value = 44250
$ws.Cells.Item(131, 4).Value = 44250
$ws.Cells.Item(132, 4).Value = 44509
$ws.Cells.Item(133, 4).Value = 44509
$ws.Cells.Item(134, 4).Value = 44488
$ws.Cells.Item(135, 4).Value = 44488
$ws.Cells.Item(136, 4).Value = 44341
$ws.Cells.Item(137, 4).Value = 44341
$ws.Cells.Item(138, 4).Value = 44278
$ws.Cells.Item(139, 4).Value = 44278
$ws.Cells.Item(140, 4).Value = 44322
$ws.Cells.Item(141, 4).Value = 44322
$ws.Cells.Item(142, 4).Value = 44194
$ws.Cells.Item(143, 4).Value = 44194
$ws.Cells.Item(144, 4).Value = 44434
$ws.Cells.Item(145, 4).Value = 44434
$ws.Cells.Item(146, 4).Value = 44490
$ws.Cells.Item(147, 4).Value = 44490
$ws.Cells.Item(148, 4).Value = 44427
$ws.Cells.Item(149, 4).Value = 44427
$ws.Cells.Item(150, 4).Value = 44266
$ws.Cells.Item(151, 4).Value = 44266
$ws.Cells.Item(152, 4).Value = 44264
$ws.Cells.Item(153, 4).Value = 44264
$ws.Cells.Item(154, 4).Value = 44525
$ws.Cells.Item(155, 4).Value = 44525
$ws.Cells.Item(156, 4).Value = 44327
$ws.Cells.Item(157, 4).Value = 44327
$ws.Cells.Item(158, 4).Value = 44462
$ws.Cells.Item(159, 4).Value = 44462

# --- Column O (Origen) updates for the rows whose region changed ---
$ws.Cells.Item(100, 15).Value = "Región de Ñuble"
$ws.Cells.Item(101, 15).Value = "Región de Ñuble"
$ws.Cells.Item(102, 15).Value = "Región Metropolitana"
$ws.Cells.Item(103, 15).Value = "Región Metropolitana"
$ws.Cells.Item(104, 15).Value = "Región de Ñuble"
$ws.Cells.Item(105, 15).Value = "Región de Ñuble"
$ws.Cells.Item(106, 15).Value = "Región Metropolitana"
$ws.Cells.Item(107, 15).Value = "Región Metropolitana"

# Make sure every Fecha cell in D2:D161 keeps the workbook's date format
$ws.Range("D2:D161").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- New row 160: Cilantro, Primera ---
$ws.Cells.Item(160, 1).Value = 11
$ws.Cells.Item(160, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(160, 3).Value = "Bíobío"
$ws.Cells.Item(160, 4).Value = 44512
$ws.Cells.Item(160, 5).Value = 8
$ws.Cells.Item(160, 6).Value = 100112040
$ws.Cells.Item(160, 7).Value = "Cilantro"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 200
$ws.Cells.Item(160, 11).Value = 600
$ws.Cells.Item(160, 12).Value = 700
$ws.Cells.Item(160, 13).Value = 650
$ws.Cells.Item(160, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(160, 15).Value = "Región de Ñuble"
$ws.Cells.Item(160, 16).Value = 650
$ws.Cells.Item(160, 17).Value = 1
$ws.Cells.Item(160, 18).Value = "Hortaliza"

# --- New row 161: Cilantro, Segunda ---
$ws.Cells.Item(161, 1).Value = 11
$ws.Cells.Item(161, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(161, 3).Value = "Bíobío"
$ws.Cells.Item(161, 4).Value = 44512
$ws.Cells.Item(161, 5).Value = 8
$ws.Cells.Item(161, 6).Value = 100112040
$ws.Cells.Item(161, 7).Value = "Cilantro"
$ws.Cells.Item(161, 8).Value = "Sin especificar"
$ws.Cells.Item(161, 9).Value = "Segunda"
$ws.Cells.Item(161, 10).Value = 100
$ws.Cells.Item(161, 11).Value = 500
$ws.Cells.Item(161, 12).Value = 500
$ws.Cells.Item(161, 13).Value = 500
$ws.Cells.Item(161, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(161, 15).Value = "Región de Ñuble"
$ws.Cells.Item(161, 16).Value = 500
$ws.Cells.Item(161, 17).Value = 1
$ws.Cells.Item(161, 18).Value = "Hortaliza"

# Re-apply the date number format to the two freshly written Fecha cells
$ws.Range("D160:D161").NumberFormat = "YYYY-MM-DD HH:MM:SS"
